$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (original row 11 was an exact
# duplicate of row 10) so everything below shifts up by one.
$ws.Rows.Item(11).Delete()

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$ws.Range("B9").Value = "Alvearie Team"

# The old "Contact" row (now row 10) becomes "Jurisdiction" /
# "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value filled in. Assigning the literal string "true"
# via .Value gets auto-coerced to the Boolean TRUE by Excel, so instead
# write it as a quoted-text formula and flatten it to a static value via
# copy / paste-special so the stored cell keeps an ordinary text type.
$caseSensitiveCell = $ws.Range("B14")
$caseSensitiveCell.Formula = '="true"'
$caseSensitiveCell.Copy() | Out-Null
$caseSensitiveCell.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
